$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style used by the other header cells (bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the I and J numeric values for rows 2-32
$data = @(
    @(5,7),
    @(8,9),
    @(9,9),
    @(6,7),
    @(4,6),
    @(4,6),
    @(9,9),
    @(8,8),
    @(6,9),
    @(8,9),
    @(6,7),
    @(7,7),
    @(8,9),
    @(7,8),
    @(4,6),
    @(9,9),
    @(3,6),
    @(6,8),
    @(8,9),
    @(8,8),
    @(7,8),
    @(8,8),
    @(5,6),
    @(5,7),
    @(7,7),
    @(7,7),
    @(8,8),
    @(1,5),
    @(1,4),
    @(4,6),
    @(1,2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
